$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Report")

# --- Row 2: "Them Size" test case, sample value XXL_2112 -> XXL_2795 ---
$ws.Range("B2").Value = "Thêm Size Mới"
$ws.Range("C2").Value = "1. Dashboard -> Menu Quản lý SP -> Kích cỡ`n2. Click Thêm Size`n3. Nhập 'XXL_2795'`n4. Lưu"
$ws.Range("D2").Value = "Input: XXL_2795"
$ws.Range("F2").Value = "Tìm thấy size 'XXL_2795': true"

# --- Row 3: "Sua ten Size" test case, sample value XXL_2112 -> XXL_2795 ---
$ws.Range("C3").Value = "1. Tìm size 'XXL_2795'`n2. Click Sửa`n3. Đổi thành 'XXL_Updated_2795'`n4. Lưu"
$ws.Range("D3").Value = "Old: XXL_2795 -> New: XXL_Updated_2795"
$ws.Range("E3").Value = "Tên size được cập nhật thành công (Tên cũ biến mất)"

# --- Row 4: "Xoa Size" test case, sample value XXL_Updated_2112 -> XXL_Updated_2795 ---
$ws.Range("C4").Value = "1. Tìm size 'XXL_Updated_2795'`n2. Click Xóa`n3. Xác nhận Alert`n4. Kiểm tra biến mất"
$ws.Range("D4").Value = "Target: XXL_Updated_2795"
$ws.Range("E4").Value = "Size bị xóa khỏi danh sách (Không còn tồn tại)"
$ws.Range("F4").Value = "Vẫn tìm thấy size 'XXL_Updated_2795': false"

# Re-fit row heights so they keep matching the sheet's default (writing the
# new multi-line strings above would otherwise leave an explicit custom
# row height behind).
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

# --- Column width changes (col B: 17.95703125 -> 14.25 ; col E: 32.22265625 -> 48.390625) ---
$ws.Columns.Item(2).ColumnWidth = 13.333333333333334
$ws.Columns.Item(5).ColumnWidth = 47.5
